$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9139728546142578
$ws.Range("B1").Value = 2.543910026550293
$ws.Range("C1").Value = 4.799402713775635
$ws.Range("D1").Value = 1.324440836906433
$ws.Range("E1").Value = 1.33961296081543
